$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# 1. Delete the row for account 005645211 / AGUINALDO / 150898.09 (originally row 2)
$ws.Rows.Item(2).Delete()

# 2. Update account 004383190 / MAFALDA balance from 31959.2 to 35000 (now row 2)
$ws.Cells.Item(2, 3).Value = 35000

# 3. Replace account 004526450 / MSD / 19500 with 004996634 / HIROKO / 11520.07 (now row 5)
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004996634"
$ws.Cells.Item(5, 2).Value = "HIROKO"
$ws.Cells.Item(5, 3).Value = 11520.07

# 4. Delete the now-duplicate row for account 004996634 / HIROKO / 858.13 (now row 26)
$ws.Rows.Item(26).Delete()
